$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3602.4  # was 4026.5293
$ws.Range("J40").Value = 2218.3  # was 2655.1428
$ws.Range("L40").Value = 2218.3  # was 2655.1428
$ws.Range("N40").Value = -2568.3  # was -3005.1428
$ws.Range("H69").Value = 7395.5  # was 7624.9473
$ws.Range("I69").Value = 8000.5  # was 8384.154
$ws.Range("J69").Value = 6548.5  # was 5980
$ws.Range("K69").Value = 24001.5  # was 25152.462
$ws.Range("L69").Value = 19645.5  # was 17940
$ws.Range("M69").Value = -23127.5  # was -24278.462
$ws.Range("N69").Value = -21393.5  # was -19688
$ws.Range("H72").Value = 7395.5  # was 7624.9473
$ws.Range("I72").Value = 8000.5  # was 8384.154
$ws.Range("J72").Value = 6548.5  # was 5980
$ws.Range("K72").Value = 72004.5  # was 75457.386
$ws.Range("L72").Value = 58936.5  # was 53820
$ws.Range("M72").Value = -67636.5  # was -71089.386
$ws.Range("N72").Value = -67672.5  # was -62556
$ws.Range("H116").Value = 9076.147999999999  # was 8705.311
$ws.Range("I116").Value = 12030.8125  # was 12030.875
$ws.Range("J116").Value = 4778.4546  # was 4612.3076
$ws.Range("K116").Value = 12030.8125  # was 12030.875
$ws.Range("L116").Value = 4778.4546  # was 4612.3076
$ws.Range("M116").Value = -8588.8125  # was -8588.875
$ws.Range("N116").Value = -11662.4546  # was -11496.3076
$ws.Range("H132").Value = 1083.6418  # was 1017.4722
$ws.Range("I132").Value = 1041.1637  # was 965.3
$ws.Range("K132").Value = 3123.4911  # was 2895.9
$ws.Range("M132").Value = -593.4911000000002  # was -365.8999999999996
$ws.Range("H138").Value = 1952.3334  # was 2634.32
$ws.Range("I138").Value = 1444.8649  # was 1630.9302
$ws.Range("J138").Value = 2255.1775  # was 3391.2632
$ws.Range("K138").Value = 4334.5947  # was 4892.7906
$ws.Range("L138").Value = 6765.532499999999  # was 10173.7896
$ws.Range("M138").Value = 805.4053000000004  # was 247.2093999999997
$ws.Range("N138").Value = -17045.5325  # was -20453.7896

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13501.345  # was 13477.884
$ws.Range("I32").Value = 12725.155  # was 12686.842
$ws.Range("K32").Value = 12725.155  # was 12686.842
$ws.Range("M32").Value = -12438.155  # was -12399.842
$ws.Range("H61").Value = 4501.375  # was 4864.643
$ws.Range("I61").Value = 1232.2  # was 1050.625
$ws.Range("K61").Value = 1232.2  # was 1050.625
$ws.Range("M61").Value = -1020.2  # was -838.625
$ws.Range("H110").Value = 9275.875  # was 5330.4
$ws.Range("I110").Value = 8633  # was 3164.8
$ws.Range("K110").Value = 8633  # was 3164.8
$ws.Range("M110").Value = -6588  # was -1119.8
$ws.Range("H132").Value = 5063.0645  # was 5194.1
$ws.Range("I132").Value = 3348.5  # was 3424.9312
$ws.Range("K132").Value = 10045.5  # was 10274.7936
$ws.Range("M132").Value = -7515.5  # was -7744.793600000001
$ws.Range("H136").Value = 4501.375  # was 4864.643
$ws.Range("I136").Value = 1232.2  # was 1050.625
$ws.Range("K136").Value = 3696.6  # was 3151.875
$ws.Range("M136").Value = -1146.6  # was -601.875

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 818  # was 938.05554
$ws.Range("I80").Value = 1049.5  # was 1157.7273
$ws.Range("J80").Value = 421.14285  # was 592.8570999999999
$ws.Range("K80").Value = 1049.5  # was 1157.7273
$ws.Range("L80").Value = 421.14285  # was 592.8570999999999
$ws.Range("M80").Value = -51.5  # was -159.7273
$ws.Range("N80").Value = -2417.14285  # was -2588.8571
$ws.Range("H83").Value = 818  # was 938.05554
$ws.Range("I83").Value = 1049.5  # was 1157.7273
$ws.Range("J83").Value = 421.14285  # was 592.8570999999999
$ws.Range("K83").Value = 5247.5  # was 5788.636500000001
$ws.Range("L83").Value = 2105.71425  # was 2964.2855
$ws.Range("M83").Value = -255.5  # was -796.6365000000005
$ws.Range("N83").Value = -12089.71425  # was -12948.2855
$ws.Range("H86").Value = 1631.0454  # was 1590.5652
$ws.Range("I86").Value = 1388.6842  # was 1354.25
$ws.Range("K86").Value = 1388.6842  # was 1354.25
$ws.Range("M86").Value = -265.6841999999999  # was -231.25
$ws.Range("H89").Value = 1631.0454  # was 1590.5652
$ws.Range("I89").Value = 1388.6842  # was 1354.25
$ws.Range("K89").Value = 6943.420999999999  # was 6771.25
$ws.Range("M89").Value = -1327.420999999999  # was -1155.25
$ws.Range("H105").Value = 1709.9524  # was 1642.826
$ws.Range("I105").Value = 1627.1765  # was 1554.6316
$ws.Range("K105").Value = 1627.1765  # was 1554.6316
$ws.Range("M105").Value = 119.8235  # was 192.3684000000001
$ws.Range("H130").Value = 70780  # was 0
$ws.Range("J130").Value = 70780  # was 0
$ws.Range("L130").Value = 70780  # was 0
$ws.Range("N130").Value = -80820  # was None

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2078.3076  # was 1913.3125
$ws.Range("I16").Value = 1915.625  # was 1882.2
$ws.Range("J16").Value = 2338.6  # was 1965.1666
$ws.Range("K16").Value = 1915.625  # was 1882.2
$ws.Range("L16").Value = 2338.6  # was 1965.1666
$ws.Range("M16").Value = -1628.625  # was -1595.2
$ws.Range("N16").Value = -2912.6  # was -2539.1666
$ws.Range("H22").Value = 546.86206  # was 527.129
$ws.Range("I22").Value = 444.58334  # was 428.92307
$ws.Range("K22").Value = 444.58334  # was 428.92307
$ws.Range("M22").Value = -94.58334000000002  # was -78.92307
$ws.Range("H31").Value = 56725.066  # was 65014.46
$ws.Range("I31").Value = 4483  # was 6070
$ws.Range("J31").Value = 91553.11  # was 82697.8
$ws.Range("K31").Value = 4483  # was 6070
$ws.Range("L31").Value = 91553.11  # was 82697.8
$ws.Range("M31").Value = -4188  # was -5775
$ws.Range("N31").Value = -92143.11  # was -83287.8
$ws.Range("H34").Value = 56725.066  # was 65014.46
$ws.Range("I34").Value = 4483  # was 6070
$ws.Range("J34").Value = 91553.11  # was 82697.8
$ws.Range("K34").Value = 4483  # was 6070
$ws.Range("L34").Value = 91553.11  # was 82697.8
$ws.Range("M34").Value = -4281  # was -5868
$ws.Range("N34").Value = -91957.11  # was -83101.8
$ws.Range("H100").Value = 38882.5  # was 38883
$ws.Range("J100").Value = 38882.5  # was 38883
$ws.Range("L100").Value = 38882.5  # was 38883
$ws.Range("N100").Value = -41046.5  # was -41047
$ws.Range("H107").Value = 3512.3137  # was 3567.9805
$ws.Range("I107").Value = 524.34784  # was 540.9091
$ws.Range("J107").Value = 5966.7144  # was 5864.3794
$ws.Range("K107").Value = 524.34784  # was 540.9091
$ws.Range("L107").Value = 5966.7144  # was 5864.3794
$ws.Range("M107").Value = 1395.65216  # was 1379.0909
$ws.Range("N107").Value = -9806.714400000001  # was -9704.3794
$ws.Range("H113").Value = 2078.3076  # was 1913.3125
$ws.Range("I113").Value = 1915.625  # was 1882.2
$ws.Range("J113").Value = 2338.6  # was 1965.1666
$ws.Range("K113").Value = 1915.625  # was 1882.2
$ws.Range("L113").Value = 2338.6  # was 1965.1666
$ws.Range("M113").Value = 254.375  # was 287.8
$ws.Range("N113").Value = -6678.6  # was -6305.1666
$ws.Range("H134").Value = 5757.32  # was 6188.478
$ws.Range("I134").Value = 6539.048  # was 6826.05
$ws.Range("J134").Value = 1653.25  # was 1938
$ws.Range("K134").Value = 19617.144  # was 20478.15
$ws.Range("L134").Value = 4959.75  # was 5814
$ws.Range("M134").Value = -17082.144  # was -17943.15
$ws.Range("N134").Value = -10029.75  # was -10884

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 926  # was 884
$ws.Range("I92").Value = 1014.5  # was 1050
$ws.Range("J92").Value = 867  # was 801
$ws.Range("K92").Value = 3043.5  # was 3150
$ws.Range("L92").Value = 2601  # was 2403
$ws.Range("M92").Value = -1795.5  # was -1902
$ws.Range("N92").Value = -5097  # was -4899
$ws.Range("I97").Value = 312.66666  # was 424
$ws.Range("J97").Value = 100  # was 95
$ws.Range("K97").Value = 937.9999799999999  # was 1272
$ws.Range("L97").Value = 300  # was 285
$ws.Range("M97").Value = -441.9999799999999  # was -776
$ws.Range("N97").Value = -1292  # was -1277
$ws.Range("J98").Value = 0  # was 400
$ws.Range("L98").Value = 0  # was 1200
$ws.Range("N98").ClearContents()  # was -4196

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 3974.5715  # was 4018.35
$ws.Range("J80").Value = 3862.2856  # was 3921
$ws.Range("L80").Value = 3862.2856  # was 3921
$ws.Range("N80").Value = -5858.2856  # was -5917
$ws.Range("H83").Value = 3974.5715  # was 4018.35
$ws.Range("J83").Value = 3862.2856  # was 3921
$ws.Range("L83").Value = 19311.428  # was 19605
$ws.Range("N83").Value = -29295.428  # was -29589
$ws.Range("H104").Value = 32491.5  # was 32495
$ws.Range("J104").Value = 32491.5  # was 32495
$ws.Range("L104").Value = 32491.5  # was 32495
$ws.Range("N104").Value = -39479.5  # was -39483
$ws.Range("H132").Value = 33025.117  # was 31221.555
$ws.Range("I132").Value = 36577.344  # was 35377.168
$ws.Range("J132").Value = 12422.2  # was 10443.5
$ws.Range("K132").Value = 109732.032  # was 106131.504
$ws.Range("L132").Value = 37266.60000000001  # was 31330.5
$ws.Range("M132").Value = -107202.032  # was -103601.504
$ws.Range("N132").Value = -42326.60000000001  # was -36390.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H33").Value = 15000  # was 27249.5
$ws.Range("I33").Value = 15000  # was 27249.5
$ws.Range("K33").Value = 15000  # was 27249.5
$ws.Range("M33").Value = -14710  # was -26959.5
$ws.Range("H134").Value = 54441.8  # was 55552.25
$ws.Range("J134").Value = 54476.332  # was 56714.5
$ws.Range("L134").Value = 54476.332  # was 56714.5
$ws.Range("N134").Value = -64616.332  # was -66854.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H34").Value = 0  # was 39999
$ws.Range("I34").Value = 0  # was 39999
$ws.Range("K34").Value = 0  # was 39999
$ws.Range("M34").ClearContents()  # was -39796
$ws.Range("H37").Value = 62498.75  # was 43003.668
$ws.Range("I37").Value = 0  # was 36008
$ws.Range("J37").Value = 62498.75  # was 49999.332
$ws.Range("K37").Value = 0  # was 36008
$ws.Range("L37").Value = 62498.75  # was 49999.332
$ws.Range("M37").ClearContents()  # was -35805
$ws.Range("N37").Value = -62904.75  # was -50405.332
$ws.Range("H40").Value = 0  # was 49999
$ws.Range("I40").Value = 0  # was 49999
$ws.Range("K40").Value = 0  # was 49999
$ws.Range("M40").ClearContents()  # was -49850
$ws.Range("H49").Value = 39998  # was 42998.332
$ws.Range("I49").Value = 0  # was 48999
$ws.Range("K49").Value = 0  # was 48999
$ws.Range("M49").ClearContents()  # was -48769
$ws.Range("H70").Value = 39450  # was 39899.5
$ws.Range("J70").Value = 39450  # was 39899.5
$ws.Range("L70").Value = 39450  # was 39899.5
$ws.Range("N70").Value = -40080  # was -40529.5
$ws.Range("H73").Value = 39450  # was 39899.5
$ws.Range("J73").Value = 39450  # was 39899.5
$ws.Range("L73").Value = 39450  # was 39899.5
$ws.Range("N73").Value = -41634  # was -42083.5
$ws.Range("H96").Value = 38238.785  # was 36984.965
$ws.Range("I96").Value = 78674.766  # was 73126.5
$ws.Range("J96").Value = 3194.2666  # was 3252.8667
$ws.Range("K96").Value = 78674.766  # was 73126.5
$ws.Range("L96").Value = 3194.2666  # was 3252.8667
$ws.Range("M96").Value = -77301.766  # was -71753.5
